$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "NEURO RAD"
$wb.Worksheets.Item(2).Name = "NEURO ONC"
